# Fruta / hortaliza, semanal
# Insert two new weekly rows of "Papa" data at the top of the
# "Macroferia Regional de Talca" block (rows 254-255), pushing the
# existing rows 254-259 down to 256-261.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows before row 254, shifting old rows 254:259 down to 256:261.
$ws.Rows("254:255").Insert()

# --- New row 254 ---
$ws.Range("A254").Value2 = 5
$ws.Range("B254").Value2 = "Macroferia Regional de Talca"
$ws.Range("C254").Value2 = "Maule"
$ws.Range("D254").Value2 = 44448
$ws.Range("E254").Value2 = 7
$ws.Range("F254").Value2 = 100114001
$ws.Range("G254").Value2 = "Papa"
$ws.Range("H254").Value2 = "Asterix"
$ws.Range("I254").Value2 = "1a (guarda)"
$ws.Range("J254").Value2 = 1200
$ws.Range("K254").Value2 = 9000
$ws.Range("L254").Value2 = 9000
$ws.Range("M254").Value2 = 9000
$ws.Range("N254").Value2 = "`$/saco 25 kilos"
$ws.Range("O254").Value2 = "Región de Los Lagos"
$ws.Range("P254").Value2 = 360
$ws.Range("Q254").Value2 = 25
$ws.Range("R254").Value2 = "Hortaliza"

# --- New row 255 ---
$ws.Range("A255").Value2 = 5
$ws.Range("B255").Value2 = "Macroferia Regional de Talca"
$ws.Range("C255").Value2 = "Maule"
$ws.Range("D255").Value2 = 44448
$ws.Range("E255").Value2 = 7
$ws.Range("F255").Value2 = 100114001
$ws.Range("G255").Value2 = "Papa"
$ws.Range("H255").Value2 = "Rodeo"
$ws.Range("I255").Value2 = "1a (guarda lavada)"
$ws.Range("J255").Value2 = 1200
$ws.Range("K255").Value2 = 10000
$ws.Range("L255").Value2 = 10000
$ws.Range("M255").Value2 = 10000
$ws.Range("N255").Value2 = "`$/malla 25 kilos"
$ws.Range("O255").Value2 = "Región de Los Lagos"
$ws.Range("P255").Value2 = 400
$ws.Range("Q255").Value2 = 25
$ws.Range("R255").Value2 = "Hortaliza"
